# Auto-generated edit script: updates cryptocurrency price/volume data
# in cryptos.xlsx (Coin, Link, Price, Volume(1h) columns) to match the
# Wed Sep  4 12:53:54 UTC 2024 GitHub Actions data refresh.
#
# Columns "D" (Price) hold numeric-looking text (e.g. "56.548.14", "1.00")
# that must stay stored as literal text, not get auto-converted to numbers
# by Excel. We force that by prefixing the value with a leading apostrophe
# (the standard Excel "treat as text" marker), same as how a user would
# type it in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'56.548.14"
$ws.Range("E2").Value = "  -4.27%  "

# Row 3
$ws.Range("D3").Value = "'2.382.03"
$ws.Range("E3").Value = "  -5.05%  "

# Row 4
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").Value = "'500.11"
$ws.Range("E5").Value = "  -6.84%  "

# Row 6
$ws.Range("D6").Value = "'129.11"
$ws.Range("E6").Value = "  -4.28%  "

# Row 7
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = "  -0.34%  "

# Row 8
$ws.Range("D8").Value = "'0.551"
$ws.Range("E8").Value = "  -3.57%  "

# Row 9
$ws.Range("D9").Value = "'2.408.87"
$ws.Range("E9").Value = "  -4.20%  "

# Row 10
$ws.Range("D10").Value = "'0.0957"
$ws.Range("E10").Value = "  -4.66%  "

# Row 11
$ws.Range("E11").Value = "  -1.54%  "

# Row 12
$ws.Range("D12").Value = "'0.319"
$ws.Range("E12").Value = "  -3.56%  "

# Row 13
$ws.Range("D13").Value = "'4.67"
$ws.Range("E13").Value = "  -9.96%  "

# Row 14
$ws.Range("D14").Value = "'2.810.18"
$ws.Range("E14").Value = "  -4.80%  "

# Row 15
$ws.Range("D15").Value = "'56.423.23"
$ws.Range("E15").Value = "  -4.12%  "

# Row 16
$ws.Range("D16").Value = "'21.58"
$ws.Range("E16").Value = "  -3.79%  "

# Row 17
$ws.Range("E17").Value = "  -3.61%  "

# Row 18
$ws.Range("D18").Value = "'2.346.80"
$ws.Range("E18").Value = "  -6.63%  "

# Row 19
$ws.Range("D19").Value = "'10.16"
$ws.Range("E19").Value = "  -5.16%  "

# Row 20
$ws.Range("D20").Value = "'310.07"
$ws.Range("E20").Value = "  -3.78%  "

# Row 21
$ws.Range("D21").Value = "'4.04"
$ws.Range("E21").Value = "  -5.32%  "

# Row 22
$ws.Range("D22").Value = "'6.24"
$ws.Range("E22").Value = "  -0.52%  "

# Row 23
$ws.Range("D23").Value = "'0.996"
$ws.Range("E23").Value = "  -0.39%  "

# Row 24
$ws.Range("D24").Value = "'65.23"
$ws.Range("E24").Value = "  -0.80%  "

# Row 25
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = "  -0.41%  "

# Row 26
$ws.Range("D26").Value = "'2.496.43"
$ws.Range("E26").Value = "  -4.95%  "

# Row 28
$ws.Range("E28").Value = "  -5.82%  "

# Row 29
$ws.Range("D29").Value = "'7.24"
$ws.Range("E29").Value = "  -3.01%  "

# Row 30
$ws.Range("D30").Value = "'173.71"
$ws.Range("E30").Value = "  -0.25%  "

# Row 31
$ws.Range("D31").Value = "'0.0₃0715"
$ws.Range("E31").Value = "  -6.10%  "

# Row 32
$ws.Range("D32").Value = "'1.66"
$ws.Range("E32").Value = "  -4.49%  "

# Row 33
$ws.Range("D33").Value = "'6.13"
$ws.Range("E33").Value = "  -2.54%  "

# Row 34
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.10"
$ws.Range("E34").Value = "  -7.76%  "

# Row 35
$ws.Range("B35").Value = "USDe"
$ws.Range("C35").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.06%  "

# Row 36
$ws.Range("D36").Value = "'0.995"
$ws.Range("E36").Value = "  -0.16%  "

# Row 37
$ws.Range("D37").Value = "'17.82"
$ws.Range("E37").Value = "  -1.87%  "

# Row 38
$ws.Range("D38").Value = "'1.21"
$ws.Range("E38").Value = "  -2.17%  "

# Row 39
$ws.Range("D39").Value = "'3.78"
$ws.Range("E39").Value = "  -4.28%  "

# Row 40
$ws.Range("D40").Value = "'35.85"
$ws.Range("E40").Value = "  -2.11%  "

# Row 41
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").Value = "'0.793"
$ws.Range("E41").Value = "  -3.86%  "

# Row 42
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'1.43"
$ws.Range("E42").Value = "  -6.11%  "

# Row 43
$ws.Range("D43").Value = "'131.66"
$ws.Range("E43").Value = "  -0.22%  "

# Row 44
$ws.Range("D44").Value = "'3.35"
$ws.Range("E44").Value = "  -4.18%  "

# Row 45
$ws.Range("D45").Value = "'4.86"
$ws.Range("E45").Value = "  -3.66%  "

# Row 46
$ws.Range("D46").Value = "'0.575"
$ws.Range("E46").Value = "  -3.12%  "

# Row 47
$ws.Range("D47").Value = "'254.08"
$ws.Range("E47").Value = "  -8.19%  "

# Row 48
$ws.Range("D48").Value = "'0.0899"

# Row 49
$ws.Range("D49").Value = "'0.0485"
$ws.Range("E49").Value = "  -5.17%  "

# Row 50
$ws.Range("E50").Value = "  -4.15%  "

# Row 51
$ws.Range("D51").Value = "'0.0207"
$ws.Range("E51").Value = "  -5.96%  "
